# Auto-generated script to apply scheduled market-price refresh to Yojimbo_Profits workbook
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1854196.1
$ws.Range("I15").Value = 1854196.1
$ws.Range("K15").Value = 5562588.300000001
$ws.Range("M15").Value = -5562419.300000001
$ws.Range("H28").Value = 422.17648
$ws.Range("I28").Value = 422.17648
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 422.17648
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 62.82351999999997
$ws.Range("N28").ClearContents()
$ws.Range("H87").Value = 18103.334
$ws.Range("J87").Value = 18103.334
$ws.Range("L87").Value = 18103.334
$ws.Range("N87").Value = -20599.334
$ws.Range("H90").Value = 18103.334
$ws.Range("J90").Value = 18103.334
$ws.Range("L90").Value = 54310.00199999999
$ws.Range("N90").Value = -66790.00199999999
$ws.Range("H113").Value = 4131.2856
$ws.Range("I113").Value = 2393.8
$ws.Range("K113").Value = 2393.8
$ws.Range("M113").Value = 860.1999999999998
$ws.Range("H123").Value = 40725
$ws.Range("J123").Value = 40725
$ws.Range("L123").Value = 40725
$ws.Range("N123").Value = -50525

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1046.9546
$ws.Range("I2").Value = 1051.65
$ws.Range("K2").Value = 1051.65
$ws.Range("M2").Value = -938.6500000000001
$ws.Range("H61").Value = 4934.5386
$ws.Range("I61").Value = 5103.92
$ws.Range("J61").Value = 700
$ws.Range("K61").Value = 5103.92
$ws.Range("L61").Value = 700
$ws.Range("M61").Value = -4891.92
$ws.Range("N61").Value = -1124
$ws.Range("H88").Value = 2812.5
$ws.Range("I88").Value = 1816.6666
$ws.Range("J88").Value = 3144.4443
$ws.Range("K88").Value = 1816.6666
$ws.Range("L88").Value = 3144.4443
$ws.Range("M88").Value = -1410.6666
$ws.Range("N88").Value = -3956.4443
$ws.Range("H91").Value = 2812.5
$ws.Range("I91").Value = 1816.6666
$ws.Range("J91").Value = 3144.4443
$ws.Range("K91").Value = 1816.6666
$ws.Range("L91").Value = 3144.4443
$ws.Range("M91").Value = -412.6666
$ws.Range("N91").Value = -5952.4443
$ws.Range("H116").Value = 1046.9546
$ws.Range("I116").Value = 1051.65
$ws.Range("K116").Value = 1051.65
$ws.Range("M116").Value = 1242.35
$ws.Range("H136").Value = 4934.5386
$ws.Range("I136").Value = 5103.92
$ws.Range("J136").Value = 700
$ws.Range("K136").Value = 15311.76
$ws.Range("L136").Value = 2100
$ws.Range("M136").Value = -12761.76
$ws.Range("N136").Value = -7200

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1046.9546
$ws.Range("I3").Value = 1051.65
$ws.Range("K3").Value = 1051.65
$ws.Range("M3").Value = -937.6500000000001
$ws.Range("H15").Value = 33339.168
$ws.Range("J15").Value = 33339.168
$ws.Range("L15").Value = 33339.168
$ws.Range("N15").Value = -33793.168
$ws.Range("H19").Value = 11247.363
$ws.Range("I19").Value = 500.16666
$ws.Range("J19").Value = 24144
$ws.Range("K19").Value = 500.16666
$ws.Range("L19").Value = 24144
$ws.Range("M19").Value = -327.16666
$ws.Range("N19").Value = -24490
$ws.Range("H86").Value = 2109.9583
$ws.Range("I86").Value = 2467.625
$ws.Range("J86").Value = 1394.625
$ws.Range("K86").Value = 2467.625
$ws.Range("L86").Value = 1394.625
$ws.Range("M86").Value = -1344.625
$ws.Range("N86").Value = -3640.625
$ws.Range("H89").Value = 2109.9583
$ws.Range("I89").Value = 2467.625
$ws.Range("J89").Value = 1394.625
$ws.Range("K89").Value = 12338.125
$ws.Range("L89").Value = 6973.125
$ws.Range("M89").Value = -6722.125
$ws.Range("N89").Value = -18205.125
$ws.Range("H107").Value = 769
$ws.Range("I107").Value = 670.125
$ws.Range("K107").Value = 670.125
$ws.Range("M107").Value = 1249.875
$ws.Range("H134").Value = 3556.8223
$ws.Range("I134").Value = 3825.2163
$ws.Range("J134").Value = 2315.5
$ws.Range("K134").Value = 11475.6489
$ws.Range("L134").Value = 6946.5
$ws.Range("M134").Value = -8940.6489
$ws.Range("N134").Value = -12016.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 40271
$ws.Range("J63").Value = 40271
$ws.Range("L63").Value = 40271
$ws.Range("N63").Value = -41643
$ws.Range("H66").Value = 40271
$ws.Range("J66").Value = 40271
$ws.Range("L66").Value = 120813
$ws.Range("N66").Value = -127677
$ws.Range("H70").Value = 26300
$ws.Range("J70").Value = 26300
$ws.Range("L70").Value = 26300
$ws.Range("N70").Value = -26930
$ws.Range("H73").Value = 26300
$ws.Range("J73").Value = 26300
$ws.Range("L73").Value = 26300
$ws.Range("N73").Value = -28484
$ws.Range("H132").Value = 8932.741
$ws.Range("I132").Value = 5817.7144
$ws.Range("J132").Value = 19835.334
$ws.Range("K132").Value = 17453.1432
$ws.Range("L132").Value = 59506.00199999999
$ws.Range("M132").Value = -14923.1432
$ws.Range("N132").Value = -64566.00199999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 11816.223
$ws.Range("I68").Value = 800
$ws.Range("J68").Value = 13193.25
$ws.Range("K68").Value = 2400
$ws.Range("L68").Value = 39579.75
$ws.Range("M68").Value = -1589
$ws.Range("N68").Value = -41201.75
$ws.Range("H71").Value = 11816.223
$ws.Range("I71").Value = 800
$ws.Range("J71").Value = 13193.25
$ws.Range("K71").Value = 7200
$ws.Range("L71").Value = 118739.25
$ws.Range("M71").Value = -3144
$ws.Range("N71").Value = -126851.25
$ws.Range("H132").Value = 1325.8
$ws.Range("I132").Value = 978.125
$ws.Range("J132").Value = 1557.5834
$ws.Range("K132").Value = 8803.125
$ws.Range("L132").Value = 14018.2506
$ws.Range("M132").Value = -6273.125
$ws.Range("N132").Value = -19078.2506

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4975.5835
$ws.Range("I70").Value = 4813.35
$ws.Range("J70").Value = 5178.375
$ws.Range("K70").Value = 4813.35
$ws.Range("L70").Value = 5178.375
$ws.Range("M70").Value = -4543.35
$ws.Range("N70").Value = -5718.375
$ws.Range("H73").Value = 4975.5835
$ws.Range("I73").Value = 4813.35
$ws.Range("J73").Value = 5178.375
$ws.Range("K73").Value = 4813.35
$ws.Range("L73").Value = 5178.375
$ws.Range("M73").Value = -3877.35
$ws.Range("N73").Value = -7050.375
$ws.Range("H102").Value = 1409.8462
$ws.Range("I102").Value = 1034.8
$ws.Range("K102").Value = 1034.8
$ws.Range("M102").Value = 587.2
$ws.Range("H107").Value = 602.82355
$ws.Range("J107").Value = 961
$ws.Range("L107").Value = 961
$ws.Range("N107").Value = -4801
$ws.Range("H113").Value = 1013.4545
$ws.Range("I113").Value = 1058.5
$ws.Range("J113").Value = 893.3333
$ws.Range("K113").Value = 1058.5
$ws.Range("L113").Value = 893.3333
$ws.Range("M113").Value = 1111.5
$ws.Range("N113").Value = -5233.3333
$ws.Range("H122").Value = 1805.3334
$ws.Range("I122").Value = 1531.3846
$ws.Range("J122").Value = 2517.6
$ws.Range("K122").Value = 4594.1538
$ws.Range("L122").Value = 7552.799999999999
$ws.Range("M122").Value = -2144.1538
$ws.Range("N122").Value = -12452.8
$ws.Range("H132").Value = 1912.2222
$ws.Range("I132").Value = 1532
$ws.Range("J132").Value = 3487.4285
$ws.Range("K132").Value = 4596
$ws.Range("L132").Value = 10462.2855
$ws.Range("M132").Value = -2066
$ws.Range("N132").Value = -15522.2855

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 514.4286
$ws.Range("I46").Value = 340.5
$ws.Range("J46").Value = 584
$ws.Range("K46").Value = 340.5
$ws.Range("L46").Value = 584
$ws.Range("M46").Value = -152.5
$ws.Range("N46").Value = -960
$ws.Range("H122").Value = 4860.7
$ws.Range("I122").Value = 4941.706
$ws.Range("K122").Value = 14825.118
$ws.Range("M122").Value = -12375.118

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -31240
$ws.Range("H136").Value = 6400
$ws.Range("I136").Value = 6672.222
$ws.Range("K136").Value = 20016.666
$ws.Range("M136").Value = -17466.666
